$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.03565433333333334
$ws.Range("H2").Value = 0.106963
$ws.Range("I2").Value = 0.002412342638581826
$ws.Range("J2").Value = 0.002412342638581825
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 4.173145864056334
$ws.Range("R2").Value = 37.558312776507
$ws.Range("S2").Value = 0.0007828933924254144
$ws.Range("T2").Value = 0.0007828933924254143
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.03565433333333334
$ws.Range("H3").Value = 0.106963
$ws.Range("I3").Value = 0.002412342638581826
$ws.Range("J3").Value = 0.002412342638581825
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 3.621768511095111
$ws.Range("R3").Value = 32.595916599856
$ws.Range("S3").Value = 0.0006794535174657673
$ws.Range("T3").Value = 0.000679453517465767
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.03565433333333334
$ws.Range("H4").Value = 0.106963
$ws.Range("I4").Value = 0.002412342638581826
$ws.Range("J4").Value = 0.002412342638581825
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 5.063870489153778
$ws.Range("R4").Value = 45.57483440238401
$ws.Range("S4").Value = 0.0009499957286906442
$ws.Range("T4").Value = 0.000949995728690644
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.066615333333333
$ws.Range("H5").Value = 24.199846
$ws.Range("I5").Value = 0.5457805068380079
$ws.Range("J5").Value = 0.5457805068380079
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 944.1534665790992
$ws.Range("R5").Value = 8497.381199211894
$ws.Range("S5").Value = 0.1771257306836251
$ws.Range("T5").Value = 0.1771257306836251
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.066615333333333
$ws.Range("H6").Value = 24.199846
$ws.Range("I6").Value = 0.5457805068380079
$ws.Range("J6").Value = 0.5457805068380079
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 819.4070867136392
$ws.Range("R6").Value = 7374.663780422752
$ws.Range("S6").Value = 0.1537229741763963
$ws.Range("T6").Value = 0.1537229741763963
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.066615333333333
$ws.Range("H7").Value = 24.199846
$ws.Range("I7").Value = 0.5457805068380079
$ws.Range("J7").Value = 0.5457805068380079
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 1145.67547658037
$ws.Range("R7").Value = 10311.07928922333
$ws.Range("S7").Value = 0.2149318019779865
$ws.Range("T7").Value = 0.2149318019779865
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.677692666666666
$ws.Range("H8").Value = 20.033078
$ws.Range("I8").Value = 0.4518071505234102
$ws.Range("J8").Value = 0.4518071505234102
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 781.5876200183046
$ws.Range("R8").Value = 7034.288580164742
$ws.Range("S8").Value = 0.1466279404667309
$ws.Range("T8").Value = 0.1466279404667309
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.677692666666666
$ws.Range("H9").Value = 20.033078
$ws.Range("I9").Value = 0.4518071505234102
$ws.Range("J9").Value = 0.4518071505234102
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 678.3202703805262
$ws.Range("R9").Value = 6104.882433424737
$ws.Range("S9").Value = 0.1272547078220139
$ws.Range("T9").Value = 0.1272547078220139
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.677692666666666
$ws.Range("H10").Value = 20.033078
$ws.Range("I10").Value = 0.4518071505234102
$ws.Range("J10").Value = 0.4518071505234102
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 948.4112496014116
$ws.Range("R10").Value = 8535.701246412704
$ws.Range("S10").Value = 0.1779245022346654
$ws.Range("T10").Value = 0.1779245022346654